# Generate Report for Handoff
# Update status from "In Translation" to "Ready for handoff" and refresh the
# handoff timestamps on the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# zh-cn sheet: Status + Latest Handoff Datetime
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("E2").Value = "2016-03-22 22:43:10"

# de-de sheet: Status + Latest Handoff Datetime
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("E2").Value = "2016-03-22 22:43:14"

# Overview sheet: per-language status + overall latest handoff date
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-03-22 22:43:14"
